$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2010-18")

# --- Row 5: relabel existing "CW3M C493" entry as "CW3M C???" ---
# Set A6/B6 first so the shared-string table is built in the same
# order as the target workbook (C490 before C??? before C492).
$ws.Range("A6").Value = "CW3M C490"
$ws.Range("B6").Value = "Demo_Baseline WRB 2010-18"
$ws.Range("C6").Value = "2010-18"
$ws.Range("C6").HorizontalAlignment = -4108

$ws.Range("A5").Value = "CW3M C???"
$ws.Range("B5").Value = "Demo_Baseline WRB 2010-18"

$ws.Range("A7").Value = "CW3M C492"
$ws.Range("B7").Value = "Demo_Baseline WRB 2010-18"
$ws.Range("C7").Value = "2010-18"
$ws.Range("C7").HorizontalAlignment = -4108

# --- Row 6: new class "CW3M C490" ---
$ws.Range("D6").Value = 1149.4268596666666
$ws.Range("E6").Value = 1612.6987305555554
$ws.Range("F6").Value = 14.827981777777779
$ws.Range("G6").Value = 52.671807666666659
$ws.Range("H6").Value = 5.2565411111111109
$ws.Range("I6").Value = 9.1008572222222224
$ws.Range("J6").Value = 2.782013222222222
$ws.Range("K6").Value = 616.95203977777771
$ws.Range("L6").Value = 44.78104311111111
$ws.Range("M6").Value = 1011.6116130000001
$ws.Range("N6").Value = 1167.1041938888891
$ws.Range("O6").Value = 526861.67361111112
$ws.Range("P6").Value = 286785.73958333331
$ws.Range("Q6").Value = -0.75187466666666658
$ws.Range("R6").Value = -0.00027177777777777774

$ws.Range("D6:N6").NumberFormat = "0.00"
$ws.Range("O6:P6").NumberFormat = "0"
$ws.Range("Q6").NumberFormat = "0.00"
$ws.Range("R6").NumberFormat = "0.000000"

# --- Row 7: new class "CW3M C492" ---
$ws.Range("D7").Value = 1166.2191842222223
$ws.Range("E7").Value = 1612.6987305555554
$ws.Range("F7").Value = 12.266958333333331
$ws.Range("G7").Value = 52.671807666666659
$ws.Range("H7").Value = 5.2575743333333333
$ws.Range("I7").Value = 8.9084497777777774
$ws.Range("J7").Value = 2.7827528888888886
$ws.Range("K7").Value = 597.77580088888897
$ws.Range("L7").Value = 44.291074333333334
$ws.Range("M7").Value = 1026.5127224444443
$ws.Range("N7").Value = 1185.8666722222222
$ws.Range("O7").Value = 429770.90277777775
$ws.Range("P7").Value = 286854.66319444444
$ws.Range("Q7").Value = -0.793682
$ws.Range("R7").Value = -0.0002864444444444444

$ws.Range("D7:N7").NumberFormat = "0.00"
$ws.Range("O7").NumberFormat = "0"
$ws.Range("P7").NumberFormat = "0"
$ws.Range("Q7").NumberFormat = "0.00"
$ws.Range("R7").NumberFormat = "0.000000"

# yellow highlight fill on specific row-7 cells, matching the original
# workbook's "divergent" styling for this new model run
$ws.Range("D7").Interior.Color = 65535
$ws.Range("F7").Interior.Color = 65535
$ws.Range("K7").Interior.Color = 65535
$ws.Range("M7").Interior.Color = 65535
$ws.Range("N7").Interior.Color = 65535
$ws.Range("O7").Interior.Color = 65535

$ws.Range("O7").Select()
